$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item(1)

# Corrected / new daily figures for rows 334-374 (columns A:G).
# Rows 334-372 receive revised F (AgTests) / G (AgPosit) counts (and minor
# corrections elsewhere); rows 373-374 are two brand-new daily entries
# (2021-03-12 and 2021-03-13) appended to the table.
$data = @(
        @(44228,252094,9504,1737,4784,196709,3456),
        @(44229,254826,12313,2732,4889,131098,2998),
        @(44230,256903,11223,2077,4976,102035,3359),
        @(44231,259533,11282,2630,5050,104308,2982),
        @(44232,261774,13677,2241,5135,226054,3181),
        @(44233,263326,8282,1552,5199,656902,5465),
        @(44234,264083,3784,757,5271,383236,3290),
        @(44235,265807,9410,1724,5382,291820,3661),
        @(44236,268986,13980,3179,5502,179100,3065),
        @(44237,271473,10932,2487,5629,132215,2964),
        @(44238,273904,10402,2431,5733,135593,2487),
        @(44239,276234,13991,2330,5812,291091,3303),
        @(44240,277682,7597,1448,5885,667474,4765),
        @(44241,278254,2775,572,5952,340699,2885),
        @(44242,279696,8237,1442,6063,231786,3241),
        @(44243,282864,13050,3168,6168,159909,2750),
        @(44244,285419,10848,2555,6271,127454,2794),
        @(44245,287752,10283,2333,6350,150280,2823),
        @(44246,290457,14300,2705,6424,306586,3543),
        @(44247,292143,8277,1686,6505,719205,5237),
        @(44248,292792,3059,649,6577,306444,2805),
        @(44249,294790,9869,1998,6671,222646,3452),
        @(44250,298337,15160,3547,6775,160045,2889),
        @(44251,300775,10958,2438,6859,138168,3025),
        @(44252,303420,12065,2645,6966,157444,2602),
        @(44253,306268,15731,2848,7075,320416,3339),
        @(44254,308083,8839,1815,7189,740823,5076),
        @(44255,308925,3565,842,7270,328498,2590),
        @(44256,311002,10854,2077,7388,225991,3124),
        @(44257,314359,15111,3357,7489,186388,2744),
        @(44258,317159,12624,2800,7560,165545,2433),
        @(44259,319582,11749,2423,7665,178240,2319),
        @(44260,322104,15721,2522,7739,333434,2808),
        @(44261,323390,7611,1286,7836,745187,3798),
        @(44262,323786,2747,396,7921,341014,2250),
        @(44263,325993,11222,2207,8037,229772,2530),
        @(44264,329593,16361,3600,8146,179040,2000),
        @(44265,331571,10643,1978,8244,153960,1891),
        @(44266,333872,11176,2301,8346,170926,1763),
        @(44267,336235,14496,2363,8440,320661,2137),
        @(44268,337503,7149,1268,8528,637925,2890)
    )

$startRow = 334
for ($i = 0; $i -lt $data.Length; $i++) {
    $row = $data[$i]
    $r = $startRow + $i
    for ($j = 0; $j -lt $row.Length; $j++) {
        $ws.Cells.Item($r, $j + 1).Value2 = $row[$j]
    }
}

# The two appended rows need the same date number format (style index 2,
# numFmt "yyyy-mm-dd") that column A uses throughout the table.
$ws.Range("A373:A374").NumberFormat = $ws.Range("A372").NumberFormat
